$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-BC($row, $value) {
    $ws.Range("B$row").Value = $value
    $ws.Range("C$row").Value = $value
}

# Some of the replacement values look like a number ("1") or a date
# ("01/01/2022"); a plain .Value assignment would let Excel auto-convert
# them into a real number / date serial, whereas the source file stores
# them as plain text (shared-string) cells. Force text storage while
# preserving the existing cell look (wrap text, top-aligned, same font
# colour) instead of picking up a brand new "Text" number format.
function Set-BC-Text($row, $value) {
    $b = $ws.Range("B$row")
    $b.NumberFormat = "@"
    $b.Value = $value
    $b.Style = "Normal"
    $b.WrapText = $true
    $b.VerticalAlignment = -4160

    $c = $ws.Range("C$row")
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
    $c.WrapText = $true
    $c.VerticalAlignment = -4160
    $c.Font.Color = 255
}

# Name: Environmental Impacts and adaptation -> Environmental Impacts Assessment
Set-BC 4 "Environmental Impacts Assessment"

# Créditos-trabalho: 0 -> 1
Set-BC-Text 6 "1"

# Carga horária: 60 h -> 90 h
Set-BC 7 "90 h"

# Ativação: 01/01/2018 -> 01/01/2022
Set-BC-Text 8 "01/01/2022"

# Objetivos:
Set-BC 10 "Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental."

# Objectives:
Set-BC 11 "Provide knowledge on the fundamentals, objectives and methods of environmental impact assessment."

# Programa resumido:
Set-BC 14 "Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais."

# Short syllabus: (unchanged per diff - not modified)

# Programa:
Set-BC 16 "Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso."

# Syllabus:
Set-BC 17 "Basic concepts and definitions. Origin and dissemination of the Environmental Impact Assessment (EIA). EIA and licensing: objectives and fundamentals. Brazilian legal and institutional framework for EIA. The EIA process and its components. Steps in planning and execution of an Environmental Impact Statement. Technological and localization alternatives. Baseline studies. Impact identification and prediction techniques. Methods and criteria for determining impact significance. Environmental management plan: mitigation, compensation, monitoring. Decision making and follow-up. Case studies."

# Método:
Set-BC 19 "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

# Critério:
Set-BC 20 "Média ponderada de atividades e provas."

# Norma de recuperação:
Set-BC 21 "1 (uma) prova escrita"

# Bibliografia:
Set-BC 22 "Bibliografia básicaSÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos. Oficina de textos: São Paulo, 2013. 583p.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão. Elsevier: Rio de Janeiro, 2019. 685p.Bibliografia complementar:COMPANHIA AMBIENTAL DO ESTADO DE SÃO PAULO - CETESB. MANUAL PARA ELABORAÇÃO DE ESTUDOS PARA O LICENCIAMENTO COM AVALIAÇÃO DE IMPACTO AMBIENTAL. Departamento de Desenvolvimento de Ações Estratégicas para o Licenciamento da Diretoria I- ID - CETESB. Anexo único, 2014. 250p."
